# Refresh Coop toilet-paper crawl data: shift/replace several product rows with
# newly scraped listings and bump every row timestamp to 2022-09-07 21:01:38.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = "2022-09-07 21:01:38"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "6283679"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "Oecoplan Toilettenpapier Camomille weiss 4-lagig 6 Rollen"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/toilettenpapier/oecoplan-toilettenpapier-camomille-weiss-4-lagig-6-rollen/p/6283679"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "6Rol"
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 4
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "4.50"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "0.75/1Rol"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "Preis pro 1 Rolle"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "0.75"
$ws.Range("L3").NumberFormat = "@"
$ws.Range("L3").Value = "1Rol"
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'toilettenpapier']"
$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = "Oecoplan Toilettenpapier Camomille weiss 4-lagig 6 Rollen 4.50 Schweizer Franken"
$ws.Range("O3").NumberFormat = "@"
$ws.Range("O3").Value = "2022-09-07 21:01:38"

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "6695141"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "Prix Garantie feuchtes Toilettenpapier 2x70 Stück"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/prix-garantie-feuchtes-toilettenpapier-2x70-stueck/p/6695141"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "140ST"
$ws.Range("E4").Value = 5
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "Coop"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "2.50"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "0.02/1ST"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "0.02"
$ws.Range("N4").NumberFormat = "@"
$ws.Range("N4").Value = "Prix Garantie feuchtes Toilettenpapier 2x70 Stück 2.50 Schweizer Franken"
$ws.Range("O4").NumberFormat = "@"
$ws.Range("O4").Value = "2022-09-07 21:01:38"

# Row 5
$ws.Range("O5").NumberFormat = "@"
$ws.Range("O5").Value = "2022-09-07 21:01:38"

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "6568452"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "Super Soft Premium Mandel feucht"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/super-soft-premium-mandel-feucht/p/6568452"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "50ST"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 3.5
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "Super Soft"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "2.95"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "0.06/1ST"
$ws.Range("K6").NumberFormat = "@"
$ws.Range("K6").Value = "0.06"
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N6").NumberFormat = "@"
$ws.Range("N6").Value = "Super Soft Premium Mandel feucht 2.95 Schweizer Franken"
$ws.Range("O6").NumberFormat = "@"
$ws.Range("O6").Value = "2022-09-07 21:01:38"

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "4947421"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "Oecoplan Taschentuch Calendula Box"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/oecoplan-taschentuch-calendula-box/p/4947421"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "80ST"
$ws.Range("E7").Value = 17
$ws.Range("F7").Value = 4
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "Coop"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "2.30"
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "0.03/1ST"
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = "Preis pro 1 Stück"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "0.03"
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1ST"
$ws.Range("M7").NumberFormat = "@"
$ws.Range("M7").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = "Oecoplan Taschentuch Calendula Box 2.30 Schweizer Franken"
$ws.Range("O7").NumberFormat = "@"
$ws.Range("O7").Value = "2022-09-07 21:01:38"

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "6834305"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "Zewa Wisch&amp;Weg Haushaltspapier weiss 4 Rollen"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "/de/haushalt-tier/toiletten-haushaltpapier/haushaltspapier/zewa-wisch-weg-haushaltspapier-weiss-4-rollen/p/6834305"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "192BLT"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 5
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "Zewa"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "5.50"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("M8").NumberFormat = "@"
$ws.Range("M8").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'haushaltspapier']"
$ws.Range("N8").NumberFormat = "@"
$ws.Range("N8").Value = "Zewa Wisch&amp;Weg Haushaltspapier weiss 4 Rollen 5.50 Schweizer Franken"
$ws.Range("O8").NumberFormat = "@"
$ws.Range("O8").Value = "2022-09-07 21:01:38"

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "6724076"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "Oecoplan feuchtes Toilettenpapier Duckies natural 40 Stück"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "/de/haushalt-tier/toiletten-haushaltpapier/toilettenpapier/feuchttuecher/oecoplan-feuchtes-toilettenpapier-duckies-natural-40-stueck/p/6724076"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40ST"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 4.5
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "Duckies"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "2.95"
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "0.07/1ST"
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "Preis pro 1 Stück"
$ws.Range("K9").NumberFormat = "@"
$ws.Range("K9").Value = "0.07"
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1ST"
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'toilettenpapier', 'feuchttuecher']"
$ws.Range("N9").NumberFormat = "@"
$ws.Range("N9").Value = "Oecoplan feuchtes Toilettenpapier Duckies natural 40 Stück 2.95 Schweizer Franken"
$ws.Range("O9").NumberFormat = "@"
$ws.Range("O9").Value = "2022-09-07 21:01:38"

# Row 10
$ws.Range("O10").NumberFormat = "@"
$ws.Range("O10").Value = "2022-09-07 21:01:38"

# Row 11
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("O11").NumberFormat = "@"
$ws.Range("O11").Value = "2022-09-07 21:01:38"

# Row 12
$ws.Range("E12").Value = ""
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("O12").NumberFormat = "@"
$ws.Range("O12").Value = "2022-09-07 21:01:38"

# Row 13
$ws.Range("O13").NumberFormat = "@"
$ws.Range("O13").Value = "2022-09-07 21:01:38"

# Row 14
$ws.Range("O14").NumberFormat = "@"
$ws.Range("O14").Value = "2022-09-07 21:01:38"

# Row 15
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tempo-taschentuecher-plus-aloe-kamille-12x9-stueck/p/3180824"
$ws.Range("M15").NumberFormat = "@"
$ws.Range("M15").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Range("O15").NumberFormat = "@"
$ws.Range("O15").Value = "2022-09-07 21:01:38"

# Row 16
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("O16").NumberFormat = "@"
$ws.Range("O16").Value = "2022-09-07 21:01:38"

# Row 17
$ws.Range("I17").Value = ""
$ws.Range("J17").Value = ""
$ws.Range("K17").Value = ""
$ws.Range("L17").Value = ""
$ws.Range("O17").NumberFormat = "@"
$ws.Range("O17").Value = "2022-09-07 21:01:38"

# Row 18
$ws.Range("O18").NumberFormat = "@"
$ws.Range("O18").Value = "2022-09-07 21:01:38"

# Row 19
$ws.Range("O19").NumberFormat = "@"
$ws.Range("O19").Value = "2022-09-07 21:01:38"

# Row 20
$ws.Range("O20").NumberFormat = "@"
$ws.Range("O20").Value = "2022-09-07 21:01:38"

# Row 21
$ws.Range("O21").NumberFormat = "@"
$ws.Range("O21").Value = "2022-09-07 21:01:38"

# Row 22
$ws.Range("O22").NumberFormat = "@"
$ws.Range("O22").Value = "2022-09-07 21:01:38"

# Row 23
$ws.Range("O23").NumberFormat = "@"
$ws.Range("O23").Value = "2022-09-07 21:01:38"

# Row 24
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/taschentuecher-strong-10x10-stueck/p/4687972"
$ws.Range("M24").NumberFormat = "@"
$ws.Range("M24").Value = "['haushalt-tier', 'toiletten-haushaltpapier', 'papiertaschentuecher', 'taschentuecher']"
$ws.Range("O24").NumberFormat = "@"
$ws.Range("O24").Value = "2022-09-07 21:01:38"

# Row 25
$ws.Range("E25").Value = ""
$ws.Range("O25").NumberFormat = "@"
$ws.Range("O25").Value = "2022-09-07 21:01:38"

# Row 26
$ws.Range("O26").NumberFormat = "@"
$ws.Range("O26").Value = "2022-09-07 21:01:38"

# Row 27
$ws.Range("I27").Value = ""
$ws.Range("J27").Value = ""
$ws.Range("K27").Value = ""
$ws.Range("L27").Value = ""
$ws.Range("O27").NumberFormat = "@"
$ws.Range("O27").Value = "2022-09-07 21:01:38"

# Row 28
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = "3874909"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Oecoplan Papiertaschentücher Special-Edition Calendula 30x10 Stück"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/oecoplan-papiertaschentuecher-special-edition-calendula-30x10-stueck/p/3874909"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "30ST"
$ws.Range("F28").Value = 5
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "Coop"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "3.65"
$ws.Range("I28").NumberFormat = "@"
$ws.Range("I28").Value = "0.12/1ST"
$ws.Range("K28").NumberFormat = "@"
$ws.Range("K28").Value = "0.12"
$ws.Range("N28").NumberFormat = "@"
$ws.Range("N28").Value = "Oecoplan Papiertaschentücher Special-Edition Calendula 30x10 Stück 20% Aktion 3.65 Schweizer Franken statt 4.60 Schweizer Franken"
$ws.Range("O28").NumberFormat = "@"
$ws.Range("O28").Value = "2022-09-07 21:01:38"

# Row 29
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = "6868354"
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Tempo Bamboo Eco"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "/de/haushalt-tier/toiletten-haushaltpapier/papiertaschentuecher/taschentuecher/tempo-bamboo-eco/p/6868354"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12ST"
$ws.Range("F29").Value = 3
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "Tempo"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "3.95"
$ws.Range("I29").NumberFormat = "@"
$ws.Range("I29").Value = "0.33/1ST"
$ws.Range("K29").NumberFormat = "@"
$ws.Range("K29").Value = "0.33"
$ws.Range("N29").NumberFormat = "@"
$ws.Range("N29").Value = "Tempo Bamboo Eco 3.95 Schweizer Franken"
$ws.Range("O29").NumberFormat = "@"
$ws.Range("O29").Value = "2022-09-07 21:01:38"

# Row 30
$ws.Range("E30").Value = ""
$ws.Range("O30").NumberFormat = "@"
$ws.Range("O30").Value = "2022-09-07 21:01:38"

# Row 31
$ws.Range("E31").Value = ""
$ws.Range("O31").NumberFormat = "@"
$ws.Range("O31").Value = "2022-09-07 21:01:38"

# Row 32
$ws.Range("O32").NumberFormat = "@"
$ws.Range("O32").Value = "2022-09-07 21:01:38"

# Row 33
$ws.Range("O33").NumberFormat = "@"
$ws.Range("O33").Value = "2022-09-07 21:01:38"
